$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated interpolation results (Y_UTM / X_UTM) and ZoneNumber for area 2
# rows: row index, new B (Y_UTM), new C (X_UTM). ZoneNumber (D) becomes 51 for all these rows.
$data = @(
    @(2, 141348099.896921, -312922610.203993),
    @(3, 178703408.6955554, -362295201.6112791),
    @(4, 157319701.6040834, -322206528.3678634),
    @(5, 142056384.3197211, -315497640.8960232),
    @(6, 141904992.4841089, -315481395.9970667),
    @(7, 164413949.5216025, -347491935.8960524),
    @(8, 142093030.5559932, -315170365.2597262),
    @(9, 157915916.6761881, -324720174.3036515),
    @(10, 138963431.6057103, -311358900.6426019),
    @(11, 158155895.589951, -323250782.3943694),
    @(12, 186807109.2753962, -360587711.855448),
    @(13, 171799814.8913856, -348508131.1287501),
    @(14, 137522281.9252852, -306526338.6425698),
    @(15, 148379064.2916293, -297192877.2481689),
    @(16, 187641498.9218537, -359974207.1259955),
    @(17, 175465203.9267208, -352113034.5898279),
    @(18, 181427741.6017235, -365943608.0528963),
    @(19, 189349594.7699511, -356528717.3753525),
    @(20, 141458787.0723889, -313257500.772146),
    @(21, 191972331.1404673, -361328626.4466391),
    @(22, 182904035.3541456, -348238819.2016164),
    @(23, 141382299.2509747, -314339671.5201009),
    @(24, 141721499.0632483, -314922671.6170502),
    @(25, 182286786.5278214, -355682675.3674026),
    @(26, 169006163.3432604, -340885700.9276638),
    @(27, 141238730.7640204, -313462381.3896099),
    @(28, 172961831.9667335, -331368761.8243606),
    @(29, 149682088.5391768, -304353800.1510978),
    @(30, 184511210.0317812, -351457881.0413339),
    @(31, 182290094.6041036, -355687121.3447877),
    @(32, 141202604.9046527, -313788147.2910262),
    @(33, 179846307.3456305, -365423631.8033922),
    @(34, 141535203.7030983, -314796955.9100499),
    @(35, 188638009.6649984, -355177105.0007602),
    @(36, 187073764.5966155, -359435188.227848),
    @(37, 151100699.8415523, -307056524.2853937),
    @(38, 153302545.1236823, -308638682.613012),
    @(39, 137427790.8201395, -309720944.1721357),
    @(40, 140683036.6809032, -311783616.8171786),
    @(41, 182185480.8546159, -355420101.1418816),
    @(42, 139785958.6206442, -312116013.6802651),
    @(43, 139436924.7472599, -308877426.2758054),
    @(44, 157108380.5028565, -322048073.0520066),
    @(45, 134113871.8673896, -300402910.3628743),
    @(46, 142167772.0194329, -315835162.0295094),
    @(47, 139622332.8264687, -312937151.0000179),
    @(48, 139914049.4457763, -309887552.3178197),
    @(49, 185854151.742791, -357653567.1813228),
    @(50, 142014766.7871837, -314061916.8233049),
    @(51, 159120642.1015579, -309707474.7853568),
    @(52, 165886800.0556829, -349389803.8062887),
    @(53, 176051257.9623712, -357852654.1054109),
    @(54, 134412228.8925521, -300462479.4236699),
    @(55, 141128205.7539958, -312691525.8584673),
    @(56, 174994308.3881196, -350840163.793776),
    @(57, 133152710.8121899, -303038671.2823944),
    @(58, 167500273.6470561, -351198028.2771947),
    @(59, 152566433.763808, -307027879.2079807),
    @(60, 162353705.6233473, -342984243.4647176),
    @(61, 182900000.1607319, -375730773.3307462),
    @(62, 185852166.9742152, -357650979.020691),
    @(63, 151272705.8182151, -301506965.1118869),
    @(64, 164372378.2907965, -346920916.0419698),
    @(65, 152153190.0924732, -331452261.5825568),
    @(66, 160006688.2880322, -309625249.6685805),
    @(67, 165888299.3243921, -349392130.2866428),
    @(68, 149672339.1384026, -303119257.454601),
    @(69, 149554551.1177341, -293750538.5389981),
    @(70, 156976021.7356825, -312915203.848655),
    @(71, 176153960.2739238, -343614386.3589055),
    @(72, 141461029.2286593, -313261337.4333048),
    @(73, 182287677.5470982, -355683872.88319),
    @(74, 162324654.8337574, -335285133.9150388),
    @(75, 140083517.112945, -312582969.6641316),
    @(76, 190110920.7044386, -364012734.1181236),
    @(77, 172455005.3259168, -338253843.2190355),
    @(78, 144749138.565618, -294665484.1164553),
    @(79, 133185091.1502571, -303142564.6646604),
    @(80, 167505773.2050066, -351206488.4550836),
    @(81, 165132573.56512, -348161339.9216876),
    @(82, 166376619.799163, -350320043.4346986),
    @(83, 169686995.5234621, -327114559.5798709),
    @(84, 155881092.0373029, -304734885.7252286),
    @(85, 184685398.9962063, -355546272.1909546),
    @(86, 166779154.6414888, -351000639.1919017),
    @(87, 149122844.3289035, -298040444.5277905),
    @(88, 159556139.6523196, -338145521.0098385),
    @(89, 133547828.5310905, -303836046.3182639),
    @(90, 152283664.1886912, -330938589.2846747),
    @(91, 133188057.8410134, -303147861.3758538),
    @(92, 148270735.1169394, -314006239.1560445),
    @(93, 132326280.4560874, -301332364.7394353),
    @(94, 133186820.3860261, -303145652.0372005),
    @(95, 132393994.8824257, -301545098.7286289),
    @(96, 158370105.7030456, -317524795.5387079),
    @(97, 137142253.7138787, -308841030.0071481),
    @(98, 158762612.8352239, -326218674.1676497),
    @(99, 166557395.5412509, -350828214.6929647),
    @(100, 161864592.474834, -342494323.4409122),
    @(101, 185912412.2909386, -352339609.3305225),
    @(102, 141352485.3882477, -312930112.7904053),
    @(103, 153369349.7000767, -331885890.8858885),
    @(104, 158134275.1538757, -333627310.9737434),
    @(105, 144420944.0583858, -314299043.6368001),
    @(106, 195960828.5914616, -390970843.5513169),
    @(107, 175859555.1756639, -348190757.0122535),
    @(108, 140639105.5325498, -313833282.4795468),
    @(109, 163555015.9934095, -317730700.4146501),
    @(110, 131829548.5578723, -295935671.3722316),
    @(111, 156634251.2094057, -311991914.7384276),
    @(112, 167670165.5312662, -353066820.298147),
    @(113, 134124177.6926021, -300005557.9198745),
    @(114, 133905416.5137617, -299761365.7542637),
    @(115, 189752920.4501799, -383696858.5584973),
    @(116, 156894333.2309652, -315577009.9073716),
    @(117, 157278467.8209096, -321661930.9312663),
    @(118, 181492026.8387491, -351779561.8086059),
    @(119, 158654562.4475776, -332007778.8960634),
    @(120, 138590394.9895089, -307612979.354601),
    @(121, 133185990.8395032, -303144170.9692506),
    @(122, 132580669.0374985, -297392959.7766616),
    @(123, 133511768.317774, -298980260.4610214),
    @(124, 168895924.6174865, -334874479.5921599),
    @(125, 135980696.7712135, -302259460.6362427),
    @(126, 156055763.8399591, -319124905.5774159),
    @(127, 148894692.201322, -292750297.2836502),
    @(128, 133184779.6796461, -303142008.564218),
    @(129, 157709961.4847255, -324568981.0597339),
    @(130, 136535760.4980252, -294347993.0635681),
    @(131, 149732268.4740105, -316086388.4335951),
    @(132, 134889794.6426065, -300646322.6606014),
    @(133, 132864678.2326166, -302570221.9648896),
    @(134, 152280106.6002271, -325158468.1114321),
    @(135, 159279825.8139072, -308512725.8769764),
    @(136, 140463547.6134146, -310251056.9035797),
    @(137, 155355626.9989494, -319283793.2341586),
    @(138, 146482371.0578764, -313926944.2323142),
    @(139, 177658724.1819302, -368915994.1828201),
    @(140, 162261097.3869281, -340454572.3690091),
    @(141, 134416284.8910971, -300469559.7393534),
    @(142, 164979926.9176007, -325686176.3676955),
    @(143, 155434422.7048624, -318239233.6028578),
    @(144, 140488008.7864806, -288631722.7400994),
    @(145, 180469151.9710608, -340246412.4500048),
    @(146, 135273482.2421477, -302242974.9294235),
    @(147, 160201540.2584676, -339546457.0315033),
    @(148, 190124677.2012815, -355283189.0850319),
    @(149, 160612346.5819714, -337121517.7568866),
    @(150, 143562825.7276157, -315222775.6070763),
    @(151, 157286568.1629081, -324247018.6740065),
    @(152, 193692536.7711925, -387386853.8541706),
    @(153, 164184943.0664755, -344101784.2502127),
    @(154, 141499581.2481183, -312502448.8642993),
    @(155, 179374978.14723, -350749944.9179873),
    @(156, 138669240.9957263, -306987520.5271319),
    @(157, 172370145.7890951, -334988174.8863518),
    @(158, 158047362.4552254, -308478310.5134187),
    @(159, 136034211.1299297, -303707503.5453642),
    @(160, 190018661.6349129, -355019272.5760236),
    @(161, 175674117.4082528, -334936743.4465915),
    @(162, 159862336.8049855, -338574060.1625065),
    @(163, 172501839.5428112, -337070218.076261),
    @(164, 140978618.3775643, -311806607.7948179),
    @(165, 176504865.0907648, -335774920.3028822),
    @(166, 186579111.1927735, -356316923.1990488),
    @(167, 140871400.795211, -311478364.8432201),
    @(168, 163430863.372584, -332212175.4115949),
    @(169, 185508470.2791307, -344604991.2797332),
    @(170, 168769481.6174906, -330232150.1497709),
    @(171, 174362090.1054911, -333699826.9162512),
    @(172, 185345857.7778867, -356796154.4118779),
    @(173, 190360136.8154513, -354967632.5779494),
    @(174, 155689613.3434788, -321488368.2586542),
    @(175, 186277094.8049808, -346526615.0679132),
    @(176, 188046861.7576432, -361483242.0999983),
    @(177, 157062556.0037677, -321069006.3031799),
    @(178, 161518848.3261522, -341509384.2953113),
    @(179, 160172432.4412201, -339007995.8900697),
    @(180, 178732917.1826238, -361906707.6826982),
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = 51
}
